# Generate Report for Archive
#
# 1) Replace the status text "Ready for handoff" with "In Translation"
#    everywhere it appears (Overview, zh-cn, de-de sheets).
# 2) Shrink the "Status" related columns (Overview!E:F, zh-cn!C, de-de!C)
#    from width 17.2159881591797 to 13.4101848602295.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq [string]$cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 13.4101848602295
$overview.Columns.Item(6).ColumnWidth = 13.4101848602295

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 13.4101848602295

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 13.4101848602295
